$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has data in A:D (D holds the repeated "jl. RCM" label).
# Insert a new blank column at D, which pushes the existing D column (and its
# data) one slot to the right, i.e. into E - matching the diff where D1:D6
# ("jl. RCM") become E1:E6 and D is left empty.
$ws.Columns("D:D").Insert(-4161) | Out-Null

# Give the now-empty column D a custom width (close to column C's own custom
# width of 11.1640625 - the nearest value reachable through ColumnWidth's
# pixel-quantized setter is 10.33 "display" units, which round-trips to
# 11.1666... internally, i.e. as close as this property allows).
$ws.Columns("D:D").ColumnWidth = 10.33

# Move the selection/active cell to B3, matching the saved selection.
$ws.Range("B3").Select() | Out-Null
